$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 36.25
$ws.Range("I11").Value = 36.25
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 36.25
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 103.75
$ws.Range("H17").Value = 1955.875
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 1955.875
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 5867.625
$ws.Range("N17").Value = -6203.625
$ws.Range("H48").Value = 16983
$ws.Range("I48").Value = 14950
$ws.Range("J48").Value = 17999.5
$ws.Range("K48").Value = 44850
$ws.Range("L48").Value = 53998.5
$ws.Range("M48").Value = -44558
$ws.Range("N48").Value = -54582.5
$ws.Range("H56").Value = 16983
$ws.Range("I56").Value = 14950
$ws.Range("J56").Value = 17999.5
$ws.Range("K56").Value = 44850
$ws.Range("L56").Value = 53998.5
$ws.Range("M56").Value = -44316
$ws.Range("N56").Value = -55066.5
$ws.Range("H59").Value = 0
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("L59").ClearContents()
$ws.Range("N59").ClearContents()
$ws.Range("H68").Value = 70000
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 70000
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 70000
$ws.Range("N68").Value = -71498
$ws.Range("H71").Value = 70000
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 70000
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 210000
$ws.Range("N71").Value = -217488
$ws.Range("H107").Value = 1424.7858
$ws.Range("I107").Value = 1538.6666
$ws.Range("J107").Value = 1219.8
$ws.Range("K107").Value = 1538.6666
$ws.Range("L107").Value = 1219.8
$ws.Range("M107").Value = 381.3334
$ws.Range("N107").Value = -5059.8
$ws.Range("H132").Value = 5664.6665
$ws.Range("I132").Value = 6088.727
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 18266.181
$ws.Range("L132").Value = 3000
$ws.Range("M132").Value = -15736.181
$ws.Range("N132").Value = -8060

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 25000
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 25000
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 25000
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -25458
$ws.Range("H32").Value = 3200.8096
$ws.Range("I32").Value = 2860.9
$ws.Range("J32").Value = 9999
$ws.Range("K32").Value = 2860.9
$ws.Range("L32").Value = 9999
$ws.Range("M32").Value = -2573.9
$ws.Range("N32").Value = -10573
$ws.Range("H62").Value = 47325
$ws.Range("I62").Value = 24300
$ws.Range("J62").Value = 55000
$ws.Range("K62").Value = 24300
$ws.Range("L62").Value = 55000
$ws.Range("M62").Value = -23676
$ws.Range("N62").Value = -56248
$ws.Range("H65").Value = 47325
$ws.Range("I65").Value = 24300
$ws.Range("J65").Value = 55000
$ws.Range("K65").Value = 72900
$ws.Range("L65").Value = 165000
$ws.Range("M65").Value = -69780
$ws.Range("N65").Value = -171240
$ws.Range("H97").Value = 1081.3572
$ws.Range("I97").Value = 285.36365
$ws.Range("J97").Value = 4000
$ws.Range("K97").Value = 285.36365
$ws.Range("L97").Value = 4000
$ws.Range("M97").Value = 210.63635
$ws.Range("N97").Value = -4992
$ws.Range("H132").Value = 2044.3529
$ws.Range("I132").Value = 2044.3529
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6133.0587
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3603.0587

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2159.6155
$ws.Range("I86").Value = 2279.7273
$ws.Range("J86").Value = 1499
$ws.Range("K86").Value = 2279.7273
$ws.Range("L86").Value = 1499
$ws.Range("M86").Value = -1156.7273
$ws.Range("N86").Value = -3745
$ws.Range("H89").Value = 2159.6155
$ws.Range("I89").Value = 2279.7273
$ws.Range("J89").Value = 1499
$ws.Range("K89").Value = 11398.6365
$ws.Range("L89").Value = 7495
$ws.Range("M89").Value = -5782.636500000001
$ws.Range("N89").Value = -18727
$ws.Range("H94").Value = 560
$ws.Range("I94").Value = 426.8
$ws.Range("J94").Value = 782
$ws.Range("K94").Value = 426.8
$ws.Range("L94").Value = 782
$ws.Range("M94").Value = 24.19999999999999
$ws.Range("N94").Value = -1684
$ws.Range("H99").Value = 2153.75
$ws.Range("I99").Value = 1200
$ws.Range("J99").Value = 5015
$ws.Range("K99").Value = 1200
$ws.Range("L99").Value = 5015
$ws.Range("M99").Value = 298
$ws.Range("N99").Value = -8011
$ws.Range("H105").Value = 5083.5454
$ws.Range("I105").Value = 4666.6665
$ws.Range("J105").Value = 5239.875
$ws.Range("K105").Value = 4666.6665
$ws.Range("L105").Value = 5239.875
$ws.Range("N105").Value = -8733.875
$ws.Range("M105").Value = -2919.6665

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1005.1429
$ws.Range("I94").Value = 1925
$ws.Range("J94").Value = 637.2
$ws.Range("K94").Value = 1925
$ws.Range("L94").Value = 637.2
$ws.Range("M94").Value = -1474
$ws.Range("N94").Value = -1539.2
$ws.Range("H96").Value = 22569
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 22569
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 22569
$ws.Range("N96").Value = -28061

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 5177.8
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 5177.8
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 15533.4
$ws.Range("N52").Value = -16065.4
$ws.Range("H92").Value = 301.5
$ws.Range("I92").Value = 435.5
$ws.Range("J92").Value = 167.5
$ws.Range("K92").Value = 1306.5
$ws.Range("L92").Value = 502.5
$ws.Range("M92").Value = -58.5
$ws.Range("N92").Value = -2998.5
$ws.Range("H107").Value = 822.6667
$ws.Range("I107").Value = 809.5
$ws.Range("J107").Value = 849
$ws.Range("K107").Value = 2428.5
$ws.Range("L107").Value = 2547
$ws.Range("M107").Value = -508.5
$ws.Range("N107").Value = -6387
$ws.Range("H133").Value = 16324.75
$ws.Range("I133").Value = 16433
$ws.Range("J133").Value = 16000
$ws.Range("K133").Value = 49299
$ws.Range("L133").Value = 48000
$ws.Range("M133").Value = -44239
$ws.Range("N133").Value = -58120

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3261.923
$ws.Range("I132").Value = 2717.5454
$ws.Range("J132").Value = 6256
$ws.Range("K132").Value = 8152.6362
$ws.Range("L132").Value = 18768
$ws.Range("M132").Value = -5622.6362
$ws.Range("N132").Value = -23828
$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").ClearContents()
$ws.Range("N134").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9426.286
$ws.Range("I7").Value = 9579.833000000001
$ws.Range("J7").Value = 8505
$ws.Range("K7").Value = 9579.833000000001
$ws.Range("L7").Value = 8505
$ws.Range("M7").Value = -9467.833000000001
$ws.Range("N7").Value = -8729
$ws.Range("H93").Value = 868.3333
$ws.Range("I93").Value = 802.5
$ws.Range("J93").Value = 1000
$ws.Range("K93").Value = 802.5
$ws.Range("L93").Value = 1000
$ws.Range("N93").Value = -3496
$ws.Range("M93").Value = 445.5
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").ClearContents()
$ws.Range("N94").ClearContents()
$ws.Range("H126").Value = 9426.286
$ws.Range("I126").Value = 9579.833000000001
$ws.Range("J126").Value = 8505
$ws.Range("K126").Value = 28739.499
$ws.Range("L126").Value = 25515
$ws.Range("M126").Value = -26269.499
$ws.Range("N126").Value = -30455
$ws.Range("H132").Value = 5099.3887
$ws.Range("I132").Value = 4586
$ws.Range("J132").Value = 7666.3335
$ws.Range("K132").Value = 13758
$ws.Range("L132").Value = 22999.0005
$ws.Range("M132").Value = -11228
$ws.Range("N132").Value = -28059.0005
$ws.Range("H136").Value = 2329.5557
$ws.Range("I136").Value = 1916.7142
$ws.Range("J136").Value = 3774.5
$ws.Range("K136").Value = 5750.142599999999
$ws.Range("L136").Value = 11323.5
$ws.Range("M136").Value = -3200.142599999999
$ws.Range("N136").Value = -16423.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3625.913
$ws.Range("I81").Value = 3266
$ws.Range("J81").Value = 3752.9412
$ws.Range("K81").Value = 6532
$ws.Range("L81").Value = 7505.8824
$ws.Range("M81").Value = -5471
$ws.Range("N81").Value = -9627.8824
$ws.Range("H84").Value = 3625.913
$ws.Range("I84").Value = 3266
$ws.Range("J84").Value = 3752.9412
$ws.Range("K84").Value = 32660
$ws.Range("L84").Value = 37529.412
$ws.Range("M84").Value = -27356
$ws.Range("N84").Value = -48137.412
$ws.Range("H122").Value = 4497.826
$ws.Range("I122").Value = 4292.684
$ws.Range("J122").Value = 5472.25
$ws.Range("K122").Value = 12878.052
$ws.Range("L122").Value = 16416.75
$ws.Range("M122").Value = -10428.052
$ws.Range("N122").Value = -21316.75
$ws.Range("H132").Value = 2016.2174
$ws.Range("I132").Value = 2016.2174
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6048.6522
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -3518.6522
$ws.Range("H136").Value = 5467.4414
$ws.Range("I136").Value = 1817.2941
$ws.Range("J136").Value = 9117.588
$ws.Range("K136").Value = 5451.8823
$ws.Range("L136").Value = 27352.764
$ws.Range("M136").Value = -2901.8823
$ws.Range("N136").Value = -32452.764
